$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the three runs that make up the opening paragraph's text
# ("Para la comunidad ... En el " + "último" + " desafío ... de los ")
# into a single run by re-typing the identical, already-concatenated text
# over that span (forces the host to coalesce the runs).
# ---------------------------------------------------------------------------
$fullIntro = "Para la comunidad, hace un par de meses he iniciado un curso de programación donde nos han dejado algunos desafíos. En el último desafío me he encontrado con variadas novedades de este mundo, me impresiona lo lejos que se puede llevar la tecnología. Adicional a lo anterior he de comentar que recientemente se realizó una clase Live donde mostraron un poco como desenvolverse en el ejercicio y ¡fue increíble!, de la clase extraje 2 lecciones muy valiosas que me hicieron llegar a sorprenderme más con la versatilidad de los "

$rng = $d.Content
$rng.Find.Execute($fullIntro, $true, $false, $false, $false, $false, $true, 1, $false, $fullIntro, 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Usa length" was split as "Usa l" / "ength"; move the "l" so the
# split becomes "Usa " / "length". The word "length" appears elsewhere in the
# document (inside code snippets), so after the first replace we continue the
# second search from right where we left off instead of re-scanning the whole
# document, to make sure we touch only this specific occurrence.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Usa l", $true, $false, $false, $false, $false, $true, 1, $false, "Usa ", 2) | Out-Null
$rng.Collapse(0)
$rng.Find.Execute("ength", $true, $false, $false, $false, $false, $true, 1, $false, "length", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: append nine new paragraphs, each holding a single hyperlink, then
# one final, completely empty paragraph, right before the end of the body.
# ---------------------------------------------------------------------------
$urls = @(
  "https://developer.mozilla.org/es/docs/Web/API/KeyboardEvent",
  "http://w3.unpocodetodo.info/jsblog/eventos-de-teclado.php",
  "https://www.youtube.com/watch?v=NnrrbfOX2x8",
  "https://www.youtube.com/watch?v=03eid8Lc8V8",
  "https://www.w3schools.com/jsref/met_element_addeventlistener.asp",
  "https://www.w3schools.com/jsref/dom_obj_event.asp",
  "https://www.youtube.com/watch?v=IQchmLGDXgU",
  "https://www.youtube.com/watch?v=2oHVjLrnRmY&t=51s",
  "https://www.youtube.com/watch?v=r-w_0SU-I74&list=PLvq-jIkSeTUZ6QgYYO3MwG9EMqC-KoLXA&index=85"
)

foreach ($u in $urls) {
  $end = $d.Content
  $end.Collapse(0)
  $end.InsertParagraphAfter()

  $newPara = $d.Paragraphs.Last
  $newPara.Range.Text = "X"

  $newPara = $d.Paragraphs.Last
  $charRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)

  $d.Hyperlinks.Add($charRange, $u, $null, $null, $u) | Out-Null

  $addedHyperlink = $d.Hyperlinks.Item($d.Hyperlinks.Count)
  $addedHyperlink.Range.Style = "Hipervnculo"
}

# Final, fully empty paragraph after the last hyperlink paragraph.
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "X"

$newPara = $d.Paragraphs.Last
$charRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
$charRange.Delete()
